$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Update F2 value to the new hash string
$ws.Range("F2").Value = "5edfa2692bdacc5e6ee805c626c50cb44cebb065f092d9a1067d89f74dacd326"

# Delete row 3 entirely (shrinks used range/dimension to A1:F2)
$ws.Rows(3).Delete()
